$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-02-12 Wednesday" "2025-02-13 Thursday"

Replace-Text "964×4=" "915×4="
Replace-Text "629×5=" "786×8="
Replace-Text "208×4=" "452×8="
Replace-Text "297×2=" "705×4="
Replace-Text "159×7=" "181×4="
Replace-Text "801×9=" "126×2="
Replace-Text "435×2=" "749×8="
Replace-Text "461×2=" "686×8="
Replace-Text "993×3=" "763×4="
Replace-Text "850×2=" "704×8="
Replace-Text "895×4=" "240×5="
Replace-Text "758×9=" "978×5="
Replace-Text "305×4=" "136×3="
Replace-Text "851×2=" "356×3="
Replace-Text "167×3=" "874×3="
Replace-Text "355×8=" "972×2="
Replace-Text "371×5=" "376×9="
Replace-Text "963×7=" "646×4="
Replace-Text "134×8=" "173×9="
Replace-Text "244×3=" "954×8="
Replace-Text "879×2=" "357×9="
Replace-Text "410×9=" "998×6="
Replace-Text "677×4=" "606×6="
Replace-Text "378×6=" "696×6="
Replace-Text "830×4=" "115×9="
